$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = "28,57 TL - 28,57 TL"
$ws.Range("I2").Value = ""

# Row 3
$ws.Range("D3").Value = ""
$ws.Range("I3").Value = ""

# Row 4
$ws.Range("D4").Value = ""
$ws.Range("I4").Value = ""

# Row 5
$ws.Range("D5").Value = ""
$ws.Range("I5").Value = ""

# Row 6
$ws.Range("D6").Value = ""
$ws.Range("I6").Value = ""

# Row 7
$ws.Range("F7").Value = "%3"

# Row 8
$ws.Range("D8").Value = ""
$ws.Range("I8").Value = ""

# Row 9
$ws.Range("D9").Value = ""
$ws.Range("I9").Value = ""

# Row 10
$ws.Range("D10").Value = ""
$ws.Range("I10").Value = ""

# Row 11
$ws.Range("D11").Value = ""
$ws.Range("I11").Value = ""

# Row 12
$ws.Range("D12").Value = ""

# Row 13
$ws.Range("D13").Value = ""
$ws.Range("I13").Value = ""

# Row 14
$ws.Range("D14").Value = ""
